# Documentation updates:
#  - Rename "DigiKey PartList" -> "DigiKey Robot PartList"
#  - Make the DigiKey sheet the active tab (was "Electrical - Non Digikey")
#  - Update each sheet's remembered selection

$wb = $excel.ActiveWorkbook

# "Electrical - Non Digikey" (2nd sheet) keeps its own selection but is no
# longer the active/selected tab.
$wsElectrical = $wb.Worksheets.Item(2)
$wsElectrical.Range("E10").Select()

# "DigiKey PartList" (3rd sheet) gets renamed and becomes the active tab.
$wsDigiKey = $wb.Worksheets.Item(3)
$wsDigiKey.Name = "DigiKey Robot PartList"
$wsDigiKey.Activate()
$wsDigiKey.Range("F25").Select()
